# Updates the Niagara_A transition probability matrix (Sheet1) with
# recalculated values reflecting additional simulated games.
# (commit: "added more games, sped up simulate game logic, and drafted
# optimization logic")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 0.1561338289962825
    "C2" = 0.6282527881040892
    "J2" = 0.007434944237918215
    "P2" = 0.1226765799256506
    "S2" = 0.08550185873605948
    "C3" = 0.005813953488372093
    "J3" = 0.03488372093023256
    "P3" = 0.8255813953488372
    "S3" = 0.1337209302325581
    "P4" = 0.8095238095238095
    "S4" = 0.1904761904761905
    "B6" = 0.0823045267489712
    "D6" = 0.0205761316872428
    "F6" = 0.09876543209876543
    "J6" = 0.2880658436213992
    "O6" = 0.01234567901234568
    "Q6" = 0.1440329218106996
    "R6" = 0.0411522633744856
    "S6" = 0.3127572016460906
    "B7" = 0.1266666666666667
    "D7" = 0.02666666666666667
    "F7" = 0.04666666666666667
    "J7" = 0.08666666666666667
    "O7" = 0.01333333333333333
    "Q7" = 0.2333333333333333
    "R7" = 0.05333333333333334
    "S7" = 0.4133333333333333
    "B8" = 0.1118012422360248
    "D8" = 0.01863354037267081
    "F8" = 0.06625258799171843
    "J8" = 0.07867494824016563
    "O8" = 0.02070393374741201
    "Q8" = 0.1884057971014493
    "R8" = 0.08074534161490683
    "S8" = 0.4347826086956522
    "B9" = 0.08695652173913043
    "D9" = 0.02898550724637681
    "F9" = 0.05797101449275362
    "J9" = 0.09420289855072464
    "O9" = 0.02173913043478261
    "Q9" = 0.2391304347826087
    "R9" = 0.08695652173913043
    "S9" = 0.3840579710144927
    "B10" = 0.1151574803149606
    "D10" = 0.02066929133858268
    "E10" = 0.001968503937007874
    "F10" = 0.06791338582677166
    "J10" = 0.08661417322834646
    "O10" = 0.01771653543307087
    "Q10" = 0.2421259842519685
    "R10" = 0.08070866141732283
    "S10" = 0.3671259842519685
    "G11" = 0.1717557251908397
    "J11" = 0.07251908396946564
    "K11" = 0.2366412213740458
    "L11" = 0.5076335877862596
    "S11" = 0.01145038167938931
    "G12" = 0.7014925373134329
    "J12" = 0.2313432835820896
    "L12" = 0.02238805970149254
    "S12" = 0.04477611940298507
    "G13" = 0.6052631578947368
    "J13" = 0.3684210526315789
    "S13" = 0.02631578947368421
    "F15" = 0.04210526315789474
    "H15" = 0.1894736842105263
    "I15" = 0.07368421052631578
    "J15" = 0.3526315789473684
    "K15" = 0.06842105263157895
    "M15" = 0.02631578947368421
    "O15" = 0.07368421052631578
    "S15" = 0.1736842105263158
    "F16" = 0.03414634146341464
    "H16" = 0.2048780487804878
    "I16" = 0.08780487804878048
    "J16" = 0.4048780487804878
    "K16" = 0.0975609756097561
    "M16" = 0.01951219512195122
    "O16" = 0.04878048780487805
    "S16" = 0.1024390243902439
    "F17" = 0.04054054054054054
    "H17" = 0.2252252252252252
    "I17" = 0.07657657657657657
    "J17" = 0.4009009009009009
    "K17" = 0.09009009009009009
    "M17" = 0.002252252252252252
    "O17" = 0.07207207207207207
    "S17" = 0.09234234234234234
    "F18" = 0.03973509933774835
    "H18" = 0.2384105960264901
    "I18" = 0.06622516556291391
    "J18" = 0.3576158940397351
    "K18" = 0.08609271523178808
    "M18" = 0.01986754966887417
    "O18" = 0.09933774834437085
    "S18" = 0.09271523178807947
    "F19" = 0.03561387066541706
    "H19" = 0.2567947516401125
    "I19" = 0.05998125585754452
    "J19" = 0.3252108716026242
    "K19" = 0.1040299906279288
    "M19" = 0.02624179943767573
    "O19" = 0.05716963448922212
    "S19" = 0.1349578256794752
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
